$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update draw results: append the latest "Pick 4" draw as a new row
# right after the last existing row (row 88 -> new row 89).
#
# Columns A (Date) and C (Phase) look numeric/date-like to Excel's input
# parser ("2025-12-14", "251214"), so a leading apostrophe is used to force
# them to be stored as literal text, exactly like every other row in this
# sheet (Date/Game/Phase/Result/InsertedAt are all plain text columns).
# Columns B, D and E are not number/date-like, so they are written as-is.
$row = 89
$ws.Range("A$row").Value = "'2025-12-14"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "'251214"
$ws.Range("D$row").Value = "5-6-8-1"
$ws.Range("E$row").Value = "2025-12-14T21:38:22.349+04:00"
